# Added Raff's Alpha 1 Trace Matrix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout: A=Feature/User Story, B=Author, C=Category, D=Path,
#                E=Asset, F=Module/Class, G=Lines

$rows = @(
    @{ Row = 7;  Author = "Raffaele";  Category = "Art/Design"; Asset = "Short/long wall 2-4,  floor 2-5";                                                  Lines = "not really" },
    @{ Row = 8;  Author = "Raffaele";  Category = "Art/Design"; Asset = "Blueprinted rooms 2-6";                                                             Lines = "not really" },
    @{ Row = 9;  Author = "Raffaele";  Category = "Art/Design"; Asset = "All assets except the atomic one";                                                  Lines = "no" },
    @{ Row = 10; Author = "Raffaele";  Category = "Art/Design"; Asset = "Particle system for loot";                                                          Lines = "no" },
    @{ Row = 11; Author = "Raffaele";  Category = "Art/Design"; Asset = "Partcile system for torch";                                                         Lines = "No" },
    @{ Row = 12; Author = "Raffaele";  Category = "Art/Design"; Asset = "Bluprinted all assets";                                                             Lines = "No" },
    @{ Row = 13; Author = "Raffaele "; Category = "Art/Design"; Asset = "UI design";                                                                         Lines = "No" },
    @{ Row = 14; Author = "Raffaele";  Category = "Art/Design"; Asset = "GDD Asset section, Game Shell and Play screen as well as tweaking every other section"; Lines = "in google drive, not in code" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.Author
    $ws.Range("C$n").Value = $r.Category
    $ws.Range("E$n").Value = $r.Asset
    $ws.Range("G$n").Value = $r.Lines
}

# Row heights to accommodate the wrapped text that was entered.
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 105

# Trailing blank rows (21-23) are no longer part of the used range.
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(21).Delete()

# Selection / scroll position matching the saved view state.
[void]$ws.Range("G14").Select()
